# Disaggregation of commodity Copper
#
# 1) Rename the shared "Copper ores and concentrates" label to "Copper"
#    on every yearly sheet (cell C7).
# 2) Rotate the D/E/F (sector) values on rows 5, 7 and 8 of every yearly
#    sheet one column to the right: new D = old F, new E = old D, new F = old E.

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # --- Rename commodity label ---
    $ws.Range("C7").Value = "Copper"

    # --- Rotate D/E/F values on rows 5, 7, 8 ---
    foreach ($row in 5, 7, 8) {
        $dCell = $ws.Cells.Item($row, 4)
        $eCell = $ws.Cells.Item($row, 5)
        $fCell = $ws.Cells.Item($row, 6)

        $dVal = $dCell.Value2
        $eVal = $eCell.Value2
        $fVal = $fCell.Value2

        $dCell.Value = $fVal
        $eCell.Value = $dVal
        $fCell.Value = $eVal
    }
}
